# Add data for 2022-02-18
# Updates the "through" date from Feb 09 to Feb 10, and bumps/adds several
# neighborhood counts in the "February 2022 (through February 10)" column (B)
# as well as a handful of other scattered month columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab to reflect the new "through" date.
$ws.Name = "Through 2022-02-10"

# Update the column header text (also drives the shared string used by B1).
$ws.Range("B1").Value = "February 2022 (through February 10)"

# Row 2 - Englewood
$ws.Range("B2").Value = 3
$ws.Range("H2").Value = 1
$ws.Range("L2").Value = 2

# Row 3 - Austin
$ws.Range("D3").Value = 6
$ws.Range("L3").Value = 4

# Row 8 - South Shore
$ws.Range("L8").Value = 2

# Row 10 - Grand Crossing
$ws.Range("B10").Value = 1

# Row 15 - Garfield Park
$ws.Range("D15").Value = 5
$ws.Range("J15").Value = 1

# Row 16 - Kenwood
$ws.Range("B16").Value = 1

# Row 23 - Chatham
$ws.Range("N23").Value = 1

# Row 31 - West Ridge
$ws.Range("D31").Value = 1

# Row 32 - Avalon Park
$ws.Range("D32").Value = 1

# Row 34 - Morgan Park
$ws.Range("L34").Value = 2

# Row 36 - Portage Park
$ws.Range("J36").Value = 1

# Row 38 - South Chicago
$ws.Range("P38").Value = 1

# Row 48 - Little Village
$ws.Range("B48").Value = 1

# Row 49 - Lake View
$ws.Range("B49").Value = 1

# Row 51 - Hermosa
$ws.Range("L51").Value = 1

# Row 64 - East Village
$ws.Range("B64").Value = 1

# Row 79 - Rush & Division
$ws.Range("F79").Value = 1

# Row 83 - Ukrainian Village
$ws.Range("F83").Value = 1
$ws.Range("L83").Value = 1
